$d = $word.ActiveDocument

# Locate the final paragraph (holds the _GoBack bookmark) -- new content is
# inserted immediately before it, so it remains the last paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionRange = $lastPara.Range
$insertionRange.InsertParagraphBefore()

$target = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$newContent = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/></w:rPr><w:t>19</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/></w:rPr><w:t xml:space="preserve"> February 2018</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:i/></w:rPr><w:t>19:25</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Today was focused on creating a connection between the planets. It had a lot of problems along the way, in a form of first not </w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>displaying the prefa</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">b of cylinder at all. I tried to implement each connection inside an array, but for current implementation it wasn’t implemented. Goal was to have a connection going from one planet to another, however at the moment all connections go from one – the beginning planet. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">The prefab was initially along the y-axis, however later during testing it didn’t display in correct direction, so the mesh was rotated 90 degrees in x-axis. This displayed the mesh in correct direction in testing, but not long enough which has to be adjusted in the code. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">It was implemented using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>localScale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>, with prefab of cylinder being created of scale 0.02f in all axes. Using transform, I calculated the distance between the two points (the two planets), and calculated the midpoint from which the prefab was displayed. Tried to calculate how to stretch it along the correct path, initially trying to work in y-axis, but after many tries and testing, it wasn’t working properly, followed by x-axis with the similar result. After testing out z-axis at the end, it worked, however the value has to be multiplied by 50.0f otherwise it is not the correct size. I am not sure to that value – with the distance being calculated before – but it has to be implemented and works at the moment.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">All points go off one planet at the moment, which has to be changed to the links coming off correct planets.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p>
'@

$target.Range.InsertXML($newContent)

# Move the "_GoBack" bookmark from the (now) final paragraph into the
# "prefab was initially ..." paragraph, right after its text, matching the
# diary entry's original bookmark placement relative to the new content.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$prefabParaIndex = $d.Paragraphs.Count - 6
$prefabPara = $d.Paragraphs.Item($prefabParaIndex)
$bmPos = $prefabPara.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
